$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2017-12-31 00:00:00"

$ws.Range("O2").Value = 27549269.41
$ws.Range("P2").Value = 112503728.31
$ws.Range("Q2").Value = 86828305.90000001
$ws.Range("R2").Value = ""

$ws.Range("S2").Value = 37332535.67
$ws.Range("T2").Value = 37332535.67
$ws.Range("U2").Value = ""

$ws.Range("V2").Value = 19081541.51
$ws.Range("W2").Value = 12392250.49
$ws.Range("X2").Value = 10654.38
$ws.Range("Y2").Value = 32136484.1
$ws.Range("Z2").Value = 31996242.32
$ws.Range("AA2").Value = 4446972.91

$ws.Range("AG2").Value = 1655553.86

$ws.Range("AP2").Value = ""
$ws.Range("AQ2").Value = ""
$ws.Range("AR2").Value = ""
$ws.Range("AS2").Value = 33357626.81
$ws.Range("AT2").Value = ""
